# Update BOM for SwitchBox: adjust Screw Terminals quantity and add two new parts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Screw Terminals quantity 9 -> 8
$ws.Range("B13").Value = 8

# New row 16: MicroSD Module (written first so the shared-string table
# picks up "MicroSD Module" before "Piezo Buzzer", matching the source order)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 1
$ws.Range("D16").Value = "MicroSD Module"

# New row 15: Piezo Buzzer
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 1
$ws.Range("D15").Value = "Piezo Buzzer"

# Update the active selection to match the saved workbook view
$ws.Range("H8").Select()
